# The workbook's data rows got reshuffled: a handful of rows swapped places
# (their full content moved to a different row, cyclically), while the row
# numbers / row-level formatting stayed put. This script re-creates that
# reshuffle by rotating the cell content of each affected group of rows.
#
# For a cycle (r0, r1, r2, ..., rn-1) the NEW content of r[i] must equal the
# OLD content of r[i+1] (wrapping around), i.e. content is rotated "up" by
# one position through the cycle. We do that with a scratch row far outside
# the used range as temporary storage, full-row-width Copy operations, and
# an explicit Clear before every paste (this engine's Range.Copy only
# overwrites cells that are non-blank in the source, so the destination
# must be blanked first or stale values from a previous, wider row would
# survive the paste).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = "A"
$lastCol = "AY"
$scratchRow = 100000

function Clear-RowRange($ws, $row) {
    $ws.Range("$firstCol$row" + ":" + "$lastCol$row").Clear()
}

function Copy-RowRange($ws, $srcRow, $dstRow) {
    $src = $ws.Range("$firstCol$srcRow" + ":" + "$lastCol$srcRow")
    $dst = $ws.Range("$firstCol$dstRow" + ":" + "$lastCol$dstRow")
    $src.Copy($dst)
}

function Rotate-Rows($ws, $rows) {
    # new content of rows[i] = old content of rows[i+1], wrapping around.
    Clear-RowRange $ws $scratchRow
    Copy-RowRange $ws $rows[0] $scratchRow

    for ($i = 0; $i -lt $rows.Length - 1; $i++) {
        $dstRow = $rows[$i]
        $srcRow = $rows[$i + 1]
        Clear-RowRange $ws $dstRow
        Copy-RowRange $ws $srcRow $dstRow
    }

    $lastRow = $rows[$rows.Length - 1]
    Clear-RowRange $ws $lastRow
    Copy-RowRange $ws $scratchRow $lastRow

    Clear-RowRange $ws $scratchRow
}

Rotate-Rows $ws @(6, 7)
Rotate-Rows $ws @(12, 13, 14)
Rotate-Rows $ws @(21, 22)
Rotate-Rows $ws @(25, 26)
Rotate-Rows $ws @(32, 33)
Rotate-Rows $ws @(36, 37, 38, 39, 40)
Rotate-Rows $ws @(43, 44)
Rotate-Rows $ws @(54, 55)
Rotate-Rows $ws @(56, 59, 58, 57)

Write-Output "Row content rotated for all affected groups."
